$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C, shifting C:N to D:O
$ws.Range("C1").EntireColumn.Insert()

# Set the new header value in C1
$ws.Range("C1").Value = "HasReflectMaxStatus"

# Update selection to match the new edit focus
$ws.Range("C1").Select() | Out-Null

Write-Host "Done"
